$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to stay as Text,
# matching the original inlineStr cell type, so force text format first.
$textAddrs = @("D5", "D6", "D7", "D9", "D10", "D11", "D15", "D17", "D20", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D36", "D39", "D40", "D41", "D42", "D45", "D46", "D47", "D48", "D51")
foreach ($a in $textAddrs) {
    $ws.Range($a).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.636.31"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "2.198.86"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "257.35"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "84.90"
$ws.Range("E6").Value = "  +13.38%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "45.16"
$ws.Range("E10").Value = "  +10.47%  "
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  +4.93%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").Value = "2.530.17"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "14.32"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "2.198.67"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "0.785"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "43.581.19"
$ws.Range("E18").Value = "  +2.09%  "
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").Value = "69.84"
$ws.Range("E20").Value = "  -1.80%  "
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("E22").Value = "  +9.79%  "
$ws.Range("D23").Value = "231.56"
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("D24").Value = "9.07"
$ws.Range("E24").Value = "  -2.90%  "
$ws.Range("D26").Value = "3.56"
$ws.Range("E26").Value = "  +5.75%  "
$ws.Range("D27").Value = "10.64"
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.27"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "38.76"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("D31").Value = "173.98"
$ws.Range("E31").Value = "  +0.66%  "
$ws.Range("D32").Value = "20.38"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").Value = "0.111"
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("E37").Value = "  +5.03%  "
$ws.Range("E38").Value = "  +4.78%  "
$ws.Range("D39").Value = "12.41"
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("D40").Value = "2.87"
$ws.Range("E40").Value = "  +5.92%  "
$ws.Range("D41").Value = "2.09"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "63.14"
$ws.Range("E42").Value = "  +5.71%  "
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("E44").Value = "  +1.33%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.0978"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "8.31"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "99.95"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "1.18"
$ws.Range("E48").Value = "  +4.54%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("E50").Value = "  +8.61%  "
$ws.Range("D51").Value = "0.429"
$ws.Range("E51").Value = "  -6.30%  "
